$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New block: rows 17-18 (SCRIPT/G01P03A/um2103.ssb) ---
# Shared-string table is filled in the same column-major order the
# original workbook used (filename, then column C top-to-bottom, then
# column D top-to-bottom, then column E top-to-bottom).
$ws.Range("A17").Value = "SCRIPT/G01P03A/um2103.ssb"
$ws.Range("B17").Value = 203
$ws.Range("B18").Value = 206

$ws.Range("C17").Value = " The [CS:P]Hidden Land[CR]?"
$ws.Range("C18").Value = " I don\'t think so...?[K] No, I\'ve never\nheard of such a place."

$ws.Range("D17").Value = " [CS:P]Сокрытые Земли[CR]?"
$ws.Range("D18").Value = " Ну, не знаю?..[K] Нет, никогда\nне слышала об этом месте."

$ws.Range("E17").Value = " [CS:P]Òïëñúóúå Èåíìé[CR]?"
$ws.Range("E18").Value = " Îô, îå èîàý?..[K] Îåó, îéëïãäà\nîå òìúšàìà ïá üóïí íåòóå."

# Row heights (auto-fit result for wrapped text at these font sizes)
$ws.Rows.Item(17).RowHeight = 43.2
$ws.Rows.Item(18).RowHeight = 21.6

# Row 16 becomes the closing/separator row of its block: pick up the
# bordered "block end" formatting (styles already used elsewhere, e.g. row 14)
# via a formats-only paste so no new style/border entries are minted.
$ws.Range("A14:E14").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the cursor / selection like the saved workbook shows
$ws.Range("C15").Select()

